$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the "Byte" count for row 6 (item "Minneskort" etc.) from 82 to 92.
$ws.Range("B6").Value = 92

# Move the active selection to I17, matching the author's final click before saving.
$ws.Range("I17").Select()
